# Applies the "basePackage" refactor:
#  - Application sheet gains a new row (C5/D5) holding the "basePackage" /
#    "com.templengine" pair.
#  - Entities sheet's old header row (Domain / com.templengine, row 3) and
#    the blank row below it (row 4) are removed, shifting the Template /
#    TemplateParameter data rows up from 5/6 to 3/4 (and the trailing
#    style-only rows up accordingly).
#  - The per-entity "H" column formulas, which used to read the base
#    package from the now-removed Entities!$C$3, now read it from
#    Application!$D$5 instead.
#  - EntityProperties' B-column INDEX/MATCH formulas automatically track
#    the Entities row shift (Excel adjusts the range references when the
#    rows are deleted).

$wb = $excel.ActiveWorkbook

$app = $wb.Worksheets.Item("Application")
$ent = $wb.Worksheets.Item("Entities")

# 1. Add the new basePackage row on the Application sheet.
$app.Range("C5").Value = "basePackage"
$app.Range("D5").Value = "com.templengine"

# 2. Remove the old "Domain / com.templengine" header row (and the blank
#    row under it) from Entities; this shifts rows 5-10 up to 3-8 and
#    Excel auto-updates every formula that referenced those rows
#    (including the cross-sheet INDEX/MATCH formulas on EntityProperties).
$ent.Rows("3:4").Delete()

# 3. The H column formulas used to read "$C$3" (the now-deleted local
#    base-package cell); point them at the new Application!$D$5 cell
#    instead, same as the rest of the workbook does.
$ent.Range("H3").Formula = '=Application!$D$5&"."&G3'
$ent.Range("H4").Formula = '=Application!$D$5&"."&G4'

# 4. Restore the selections shown in the final workbook, leaving
#    "Application" as the active tab.
$ent.Activate()
$ent.Range("M3").Select()

$app.Activate()
$app.Range("C6").Select()
